# Fix percentage values in the climate data workbook:
#  1. Column C on every "län" sheet holds a KPI1 rate currently stored as a
#     raw fraction (e.g. 0.04848084661477944). Convert it to the equivalent
#     percentage number (value * 100, rounded to 1 decimal place) so it reads
#     naturally as a plain number (e.g. 4.8).
#  2. Update the column C header text (row 1) to make the unit explicit by
#     appending " (%)".

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Header in row 1 / column C carries the KPI1 label - clarify the unit.
    $header = $ws.Cells.Item(1, 3)
    $header.Value = "KPI1: Förändringstakt andel laddbara bilar (%)"

    # Convert every data row's column C fraction into a percentage number.
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count()

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $frac = $cell.Value()
        if ($frac -ne $null) {
            $pct = [Math]::Round($frac * 100, 1)
            $cell.Value = $pct
        }
    }
}
